$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.871.14'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.619.71'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.03%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E8").Value = '  -0.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0616'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.02%  '
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.843.41'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.601.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.28%  '
$ws.Range("E15").Value = '  -2.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.869.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.40%  '
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '191.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.47'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("E23").Value = '  -1.84%  '
$ws.Range("E24").Value = '  +2.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.23%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  -3.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.60%  '
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("E31").Value = '  -2.31%  '
$ws.Range("E32").Value = '  -3.71%  '
$ws.Range("E33").Value = '  -4.62%  '
$ws.Range("E34").Value = '  -1.65%  '
$ws.Range("E35").Value = '  -2.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.122.62'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.835'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.37'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.89%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.510'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.78%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0153'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.13'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.753.92'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.747'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.84%  '
$ws.Range("E44").Value = '  -5.04%  '
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.40%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '53.97'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.41%  '
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("E51").Value = '  -3.35%  '
